$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item(1)

# --- Create the three new sheets -----------------------------------------
# Creation order (drives sheetId assignment): MailReset, InvalidLoginData,
# SearchUser -- then SearchUser is moved to sit right after LoginData.
$wsMail = $wb.Worksheets.Add($null, $login)
$wsMail.Name = "MailReset"

$wsInvalid = $wb.Worksheets.Add($null, $wsMail)
$wsInvalid.Name = "InvalidLoginData"

$wsSearch = $wb.Worksheets.Add($null, $login)
$wsSearch.Name = "SearchUser"

# Re-resolve stable handles now that all the Add()/reorder calls are done.
$wsMail = $wb.Worksheets.Item("MailReset")
$wsInvalid = $wb.Worksheets.Item("InvalidLoginData")
$wsSearch = $wb.Worksheets.Item("SearchUser")

# --- Fill InvalidLoginData -------------------------------------------------
$wsInvalid.Range("A1").Value = "UserName"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "ad"
$wsInvalid.Range("B2").Value = 123456
$wsInvalid.Range("A2").Select()

# --- Fill MailReset (with mailto hyperlink) --------------------------------
$wsMail.Range("A1").Value = "email"
$wsMail.Range("A2").Value = "aswa@gmail.com"
$wsMail.Hyperlinks.Add($wsMail.Range("A2"), "mailto:aswa@gmail.com")
$wsMail.Range("A2").Select()

# --- Fill SearchUser --------------------------------------------------------
$wsSearch.Range("A1").Value = "UserName"
$wsSearch.Range("B1").Value = "Password"
$wsSearch.Range("C1").Value = "Search"
$wsSearch.Range("A2").Value = "admin"
$wsSearch.Range("B2").Value = 123456
$wsSearch.Range("C2").Value = "ljhg"

# --- Tweak LoginData (original sheet) --------------------------------------
$login.Columns.Item(1).ColumnWidth = 9.44
$login.Range("A1:B2").Select()

# SearchUser ends up as the active tab, with the cursor left on H19.
$wsSearch.Activate()
$wsSearch.Range("H19").Select()
